$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$ws.Activate()

# --- Remove the old EmailSubject / EmailBody rows ---------------------------
# Row 25 held "EmailSubject" / "Nike daily refunds Transactions Audit".
# Row 26 held "EmailBody" / "Please find attached the audit file for the
# transactions processed." (wrapped, style applied to column B).
# Clear row 25's content, then delete row 26 entirely so everything below
# shifts up by one (GenerateAuditFile -> row26, ConsoleToDateDelay -> row27).
$ws.Range("A25:C25").ClearContents()
$ws.Rows("26:26").Delete()

# --- Add the new process-completed email settings ---------------------------
$ws.Range("A29").Value = "ProcessCompletedEmailSubject"
$ws.Range("B29").Value = "Nike daily refunds Transactions Audit"

$ws.Range("A30").Value = "ProcessComepletedEmailBody"
$ws.Range("B30").Value = "Please find attached the audit file for the transactions processed."
$ws.Range("B30").WrapText = $true

# --- Add the new exception-notification email settings ----------------------
$ws.Range("A32").Value = "BusinessRuleExceptionEmailSubject"
$ws.Range("B32").Value = "Nike Daily Refunds: Business Rule Violation Occurred"

$ws.Range("A33").Value = "SystemExceptionEmailSubject"
$ws.Range("B33").Value = "Nike Daily Refunds: Exception Occurred"

# --- Restore the view/selection state ---------------------------------------
$ws.Range("B31").Select() | Out-Null
